$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename header row columns: "<name>_old" -> "<name>_FV2210" (cols A-J)
#    and "<name>_new" -> "<name>_FV2304" (cols L-U). Column K ("diff") stays.
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2210"
    $ws.Cells.Item(1, $i + 12).Value = $baseNames[$i] + "_FV2304"
}

# 2) Turn the used range A1:U88 into an Excel Table ("Table1") with a header row,
#    so the table definition picks up the freshly renamed headers.
$tableRange = $ws.Range("A1:U88")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# 3) Freeze the header row (split/freeze pane at row 2, i.e. below row 1).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
